# BBDD.xlsx schema update
#   - rename "id" key columns to "id peli" (peliculas + puntos tables)
#   - drop "actores"/"guionista" columns from the "peliculas" table header,
#     shifting "argumento"/"duracion" two columns to the left
#   - normalise the border on the "actores" table header (premios oscar cell)
#   - add a small pk/fk cross-reference box (id peli / actor, id peli / guionista)
#   - add a csv-file naming-convention note block (peliculas_raw / peliculas_api_1998,
#     datos_imdb / datos_imdb:_1998, actors / actores eli / actores )

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# A. "peliculas" table header (row 2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "id peli"

# The table used to have "actores"/"guionista" columns (H2/I2) before
# "argumento"/"duracion" (J2/K2). Those two columns are removed, so
# argumento/duracion slide left into H2/I2, and J2/K2 become empty.
$ws.Range("J2").Copy()
$ws.Range("H2").PasteSpecial(-4122)   # xlPasteFormats (drop the old green/orange fill)
$ws.Range("H2").Value = "argumento"

$ws.Range("I2").Value = "duracion"    # style already matches (plain border)

$ws.Range("J2:K2").Clear()

# ---------------------------------------------------------------------------
# B. "actores" table header (row 5) - normalise border on premios oscar cell
# ---------------------------------------------------------------------------
$ws.Range("C5").Copy()
$ws.Range("E5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E5").Value = "premios oscar"

# ---------------------------------------------------------------------------
# C. "puntos" table header (row 11)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "id peli"

# ---------------------------------------------------------------------------
# D. New pk / fk cross-reference box (H12:L14)
# ---------------------------------------------------------------------------
$ws.Range("H12:I12").Merge()
$ws.Range("H12").Value = "pk"
$ws.Range("H12:I12").HorizontalAlignment = -4108   # xlCenter

$ws.Range("A2").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H13").Value = "id peli"

$ws.Range("A2").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value = "actor"

$ws.Range("A2").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$ws.Range("K13").Value = "id peli"

$ws.Range("A2").Copy()
$ws.Range("L13").PasteSpecial(-4122)
$ws.Range("L13").Value = "guionista"

$ws.Range("A20").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("K14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# E. csv / table naming-convention notes (H20:I27)
# ---------------------------------------------------------------------------
$ws.Range("H20").Value = "peliculas_raw"
$ws.Range("I20").Value = "peliculas_api_1998"
$ws.Range("H21").Value = "datos_imdb"
$ws.Range("I21").Value = "datos_imdb:_1998"
$ws.Range("H26").Value = "actors"
$ws.Range("I26").Value = "actores eli"
$ws.Range("I27").Value = "actores "

$ws.Range("I21").Select()
